$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Project Design Phase" heading
#   - center the paragraph
#   - drop the leading run of spaces used for manual centering
#   - split "Project Design Phase" into "Project D" / "esign Phase"
#     with the _GoBack bookmark sitting between them (this mirrors
#     Word's own "last edit position" bookmark, which moves to wherever
#     text was typed most recently)
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Alignment = 1   # wdAlignParagraphCenter -> <w:jc w:val="center"/>

$rng1 = $p1.Range
$null = $rng1.Find.Execute( `
    "                                                        ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p1 = $d.Paragraphs.Item(1)
$rng1b = $p1.Range
$null = $rng1b.Find.Execute("Project D", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1b.Collapse(0)
$null = $d.Bookmarks.Add("_GoBack", $rng1b)

# ---------------------------------------------------------------------
# Paragraph 2: "Solution Architecture" heading
#   - center the paragraph
#   - drop the leading run of spaces in front of the text
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Alignment = 1

$rng2 = $p2.Range
$null = $rng2.Find.Execute( `
    "                                                         ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# Table: "Project Name" row -> rename "News App" to "InsightStream"
# ---------------------------------------------------------------------
$t = $d.Tables.Item(1)
$projectNameCell = $t.Cell(3, 2)
$null = $projectNameCell.Range.Find.Execute( `
    "News App", $true, $false, $false, $false, $false, $true, 1, $false, "InsightStream", 2)

# ---------------------------------------------------------------------
# Table: "Team Member 3" row -> fix double space in "Yukktha  R"
# ---------------------------------------------------------------------
$memberCell = $t.Cell(8, 2)
$null = $memberCell.Range.Find.Execute( `
    "Yukktha  R", $true, $false, $false, $false, $false, $true, 1, $false, "Yukktha R", 2)
